$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column price values and E-column volume-change percentages updated per latest crypto snapshot.
# Numeric-looking price strings are entered with a leading apostrophe so Excel keeps them as text
# (matching the source data's text formatting) instead of silently coercing to Number.

$ws.Range("D2").Value = "66.894.50"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "3.078.13"
$ws.Range("E3").Value = "  -1.48%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'577.99"
$ws.Range("D6").Value = "'169.33"
$ws.Range("E6").Value = "  -2.58%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.073.80"
$ws.Range("E8").Value = "  -1.44%  "
$ws.Range("D9").Value = "'0.514"
$ws.Range("E9").Value = "  -1.65%  "
$ws.Range("D10").Value = "'6.41"
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("D11").Value = "'0.149"
$ws.Range("E11").Value = "  -3.44%  "
$ws.Range("D12").Value = "'0.472"
$ws.Range("E12").Value = "  -1.62%  "
$ws.Range("D13").Value = "'0.0000241"
$ws.Range("E13").Value = "  -3.04%  "
$ws.Range("D14").Value = "'36.37"
$ws.Range("E14").Value = "  -2.54%  "
$ws.Range("E15").Value = "  -2.22%  "
$ws.Range("D16").Value = "3.593.10"
$ws.Range("D17").Value = "66.856.78"
$ws.Range("E17").Value = "  -0.37%  "
$ws.Range("D18").Value = "'7.02"
$ws.Range("E18").Value = "  -1.59%  "
$ws.Range("D19").Value = "3.081.52"
$ws.Range("E19").Value = "  -1.48%  "
$ws.Range("D20").Value = "'16.35"
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("D21").Value = "'482.04"
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("D22").Value = "'7.72"
$ws.Range("E22").Value = "  +0.51%  "
$ws.Range("D23").Value = "'0.689"
$ws.Range("E23").Value = "  -3.37%  "
$ws.Range("D24").Value = "'82.91"
$ws.Range("E24").Value = "  -1.25%  "
$ws.Range("D25").Value = "'12.84"
$ws.Range("E25").Value = "  -3.53%  "
$ws.Range("D26").Value = "'2.23"
$ws.Range("E26").Value = "  -3.59%  "
$ws.Range("D27").Value = "'10.30"
$ws.Range("E27").Value = "  +2.89%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").Value = "'7.72"
$ws.Range("E29").Value = "  -3.40%  "
$ws.Range("D30").Value = "'2.31"
$ws.Range("E30").Value = "  -3.67%  "
$ws.Range("D31").Value = "'2.63"
$ws.Range("E31").Value = "  -1.72%  "
$ws.Range("D32").Value = "'27.82"
$ws.Range("E32").Value = "  -3.17%  "
$ws.Range("D33").Value = "'0.111"
$ws.Range("E33").Value = "  -2.03%  "
$ws.Range("D34").Value = "0.0₃0916"
$ws.Range("E34").Value = "  -7.48%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("D36").Value = "'5.68"
$ws.Range("E36").Value = "  -3.15%  "
$ws.Range("D37").Value = "'0.954"
$ws.Range("E37").Value = "  -2.92%  "
$ws.Range("D38").Value = "'46.22"
$ws.Range("E38").Value = "  -3.16%  "
$ws.Range("D39").Value = "'0.123"
$ws.Range("E39").Value = "  +0.30%  "
$ws.Range("D40").Value = "'1.98"
$ws.Range("E40").Value = "  -5.03%  "
$ws.Range("D41").Value = "'0.301"
$ws.Range("E41").Value = "  -3.20%  "
$ws.Range("D42").Value = "'8.32"
$ws.Range("E42").Value = "  -3.51%  "
$ws.Range("D43").Value = "2.771.10"
$ws.Range("E43").Value = "  -2.36%  "
$ws.Range("D44").Value = "'377.35"
$ws.Range("E44").Value = "  -0.96%  "
$ws.Range("D45").Value = "'2.53"
$ws.Range("E45").Value = "  -3.92%  "
$ws.Range("D46").Value = "'0.0346"
$ws.Range("E46").Value = "  -3.13%  "
$ws.Range("D47").Value = "'135.33"
$ws.Range("E47").Value = "  -0.50%  "
$ws.Range("D49").Value = "'24.41"
$ws.Range("E49").Value = "  -1.85%  "
$ws.Range("D50").Value = "'2.15"
$ws.Range("E50").Value = "  -2.79%  "
$ws.Range("E51").Value = "  -1.94%  "
